$d = $word.ActiveDocument

# The version/date line currently reads:
#   "Version 11.02.03, 2015-05-27"
# and must become:
#   "Version 11.03.10, 2015-06-29"
#
# The original text is split across many single/few-character runs
# (one run per rsid group). We edit those runs in place, working from
# the end of the string back towards the start, because replacing the
# text of a Range only ever merges that run with runs that come AFTER
# it (never runs before it) -- so working right-to-left keeps all of
# our not-yet-processed offsets valid and keeps the already-processed
# runs intact.

$full = $d.Content.Text
$base = $full.IndexOf("Version 11.02.03, 2015-05-27")

# Append a brand new run holding "29" right after the existing "27" run.
# First insert the characters (this merges into the preceding run because
# the formatting is identical), then nudge the new span's font size away
# from and back to its original value so the engine is forced to keep it
# as a separate <w:r> with otherwise-identical run properties.
$insertAt = $base + 28
$newRun = $d.Range($insertAt, $insertAt)
$newRun.InsertAfter("29")
$splitRange = $d.Range($insertAt, $insertAt + 2)
$splitRange.Font.Size = 10
$splitRange.Font.Size = 9

# Now rewrite each old run's text, right-to-left.
$d.Range($base + 26, $base + 28).Text = "-"    # "27" -> "-"
$d.Range($base + 25, $base + 26).Text = "6"    # "-"  -> "6"
$d.Range($base + 24, $base + 25).Text = "0"    # "5"  -> "0"
$d.Range($base + 23, $base + 24).Text = "-"    # "0"  -> "-"
$d.Range($base + 22, $base + 23).Text = "5"    # "-"  -> "5"
$d.Range($base + 21, $base + 22).Text = "1"    # "5"  -> "1"
$d.Range($base + 20, $base + 21).Text = ", 20" # "1"  -> ", 20"
$d.Range($base + 16, $base + 20).Text = "0"    # ", 20" -> "0"
$d.Range($base + 15, $base + 16).Text = "1"    # "3"  -> "1"
$d.Range($base + 13, $base + 15).Text = "."    # ".0" -> "."
$d.Range($base + 12, $base + 13).Text = "3"    # "2"  -> "3"
